# Weekly CompStat (62nd Precinct) refresh: "New crime data collected"
# Updates the report header (volume number + week-covering dates) and
# refreshes every weekly/28-day/YTD/2-year crime statistic in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# Header: "Volume 30   Number  47" -> "...48"
#         "Report Covering the Week  11/20/2023  Through  11/26/2023"
#      -> "...11/27/2023  Through  12/3/2023"
# Use Characters() so the run keeps its own explicit font (the A8 cell's
# base style is 20pt, but the rich-text run is 10pt Andale WT).
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "48"
$a8 = $ws.Range("A8").Characters()
$a8.Font.Name = "Andale WT"
$a8.Font.Size = 10

$ws.Range("C9").Characters(27, 10).Text = "11/27/2023"
$ws.Range("C9").Characters(48, 10).Text = "12/3/2023"
$c9 = $ws.Range("C9").Characters()
$c9.Font.Name = "Andale WT"
$c9.Font.Size = 10

# ---------------------------------------------------------------------
# Crime-complaints table (rows 14-30): refreshed weekly/28-day/YTD/2-year
# figures. Pure numeric updates first (style/type unchanged).
# ---------------------------------------------------------------------
$ws.Range("M14").Value = -71.428571428571

$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 11.764705882352
$ws.Range("L15").Value = -5
$ws.Range("M15").Value = 72.727272727272
$ws.Range("N15").Value = -20.833333333333

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 134
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = 18.584070796460
$ws.Range("L16").Value = 34
$ws.Range("M16").Value = -16.25
$ws.Range("N16").Value = -81.100141043723

$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 69.230769230769
$ws.Range("I17").Value = 210
$ws.Range("J17").Value = 171
$ws.Range("K17").Value = 22.807017543859
$ws.Range("L17").Value = 27.272727272727
$ws.Range("M17").Value = 69.354838709677
$ws.Range("N17").Value = -27.586206896551

$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 202
$ws.Range("K18").Value = 14.772727272727
$ws.Range("L18").Value = 46.376811594202
$ws.Range("M18").Value = -24.907063197026
$ws.Range("N18").Value = -86.780104712041

$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 44.444444444444
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 5.714285714285
$ws.Range("I19").Value = 595
$ws.Range("J19").Value = 677
$ws.Range("K19").Value = -12.112259970457
$ws.Range("L19").Value = 10.594795539033
$ws.Range("M19").Value = 45.121951219512
$ws.Range("N19").Value = -18.044077134986

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 18
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 152
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = 8.571428571428
$ws.Range("L20").Value = 47.572815533980
$ws.Range("M20").Value = 2.702702702702
$ws.Range("N20").Value = -90.989922940130

$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 5.555555555555
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = 27.710843373494
$ws.Range("I21").Value = 1314
$ws.Range("J21").Value = 1299
$ws.Range("K21").Value = 1.154734411085
$ws.Range("L21").Value = 23.264540337711
$ws.Range("M21").Value = 16.386182462356
$ws.Range("N21").Value = -73.577317514578

$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 500

$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 147
$ws.Range("H24").Value = -39.455782312925
$ws.Range("I24").Value = 1565
$ws.Range("J24").Value = 1744
$ws.Range("K24").Value = -10.263761467889
$ws.Range("L24").Value = 33.418584825234
$ws.Range("M24").Value = 70.851528384279

$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 28.571428571428
$ws.Range("I25").Value = 581
$ws.Range("J25").Value = 432
$ws.Range("K25").Value = 34.490740740740
$ws.Range("L25").Value = 40.338164251207
$ws.Range("M25").Value = 40.338164251207

$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 25
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = 4.166666666666
$ws.Range("L26").Value = -10.714285714285

$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 46
$ws.Range("K27").Value = -2.127659574468
$ws.Range("L27").Value = -19.298245614035

$ws.Range("L30").Value = -73.333333333333

# ---------------------------------------------------------------------
# Cells whose type/style also changes (numeric <-> the "0" / "***.*"
# placeholder text, styled like style 14). Setting .Value alone leaves
# the old number-format style behind, so re-apply the donor cell's
# format afterwards (PasteSpecial Formats keeps the value/type intact).
# C23 carries style 14 untouched by this edit, so it's a safe donor for
# every "become text" cell; L28/L29 (style 16) are safe donors for the
# two cells that become plain numbers.
# ---------------------------------------------------------------------
function Set-TextLikeStyle14($ref, $text) {
    # A bare numeric-looking string (e.g. "0") is auto-coerced back to a
    # number by value assignment; force text with a leading quote prefix
    # so it resolves to the shared string instead, same as "***.*" (which
    # isn't numeric-looking and doesn't need the prefix).
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($ref).Value = "'" + $text
    } else {
        $ws.Range($ref).Value = $text
    }
    $ws.Range("C23").Copy()
    $ws.Range($ref).PasteSpecial($fmtFormats)
}

Set-TextLikeStyle14 "D15" "0"
Set-TextLikeStyle14 "E15" "***.*"
Set-TextLikeStyle14 "D18" "0"
Set-TextLikeStyle14 "E18" "***.*"
Set-TextLikeStyle14 "C22" "0"
Set-TextLikeStyle14 "D22" "0"
Set-TextLikeStyle14 "E22" "***.*"
Set-TextLikeStyle14 "C27" "0"
Set-TextLikeStyle14 "C30" "0"

$ws.Range("M28").Value = 100
$ws.Range("L28").Copy()
$ws.Range("M28").PasteSpecial($fmtFormats)

$ws.Range("M29").Value = 100
$ws.Range("L29").Copy()
$ws.Range("M29").PasteSpecial($fmtFormats)

$excel.CutCopyMode = $false
